$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 390
$ws.Cells.Item(390, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(390, 2).Value = 500
$ws.Cells.Item(390, 3).Value = 2900
$ws.Cells.Item(390, 4).Value = 0
$ws.Cells.Item(390, 5).Value = 0.1
$ws.Cells.Item(390, 6).Value = 47.97993275550777
$ws.Cells.Item(390, 7).Value = 1056.614746349956
$ws.Cells.Item(390, 8).Value = 28.32741093155056
$ws.Cells.Item(390, 9).Value = 9.062060965211266
$ws.Cells.Item(390, 10).Value = 28.71458278916847

# Row 391
$ws.Cells.Item(391, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(391, 2).Value = 65
$ws.Cells.Item(391, 3).Value = 2900
$ws.Cells.Item(391, 4).Value = 0
$ws.Cells.Item(391, 5).Value = 0.1
$ws.Cells.Item(391, 6).Value = 49.65609781879013
$ws.Cells.Item(391, 7).Value = 1622.896866179913
$ws.Cells.Item(391, 8).Value = 45.73203552001451
$ws.Cells.Item(391, 9).Value = 7.385895901928905
$ws.Cells.Item(391, 10).Value = 11.30995820070452

# Row 392
$ws.Cells.Item(392, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(392, 2).Value = 580
$ws.Cells.Item(392, 3).Value = 2900
$ws.Cells.Item(392, 4).Value = 0
$ws.Cells.Item(392, 5).Value = 0.1
$ws.Cells.Item(392, 6).Value = 47.11642464565518
$ws.Cells.Item(392, 7).Value = 1013.447515844847
$ws.Cells.Item(392, 8).Value = 25.25085120935917
$ws.Cells.Item(392, 9).Value = 9.925569075063855
$ws.Cells.Item(392, 10).Value = 31.79114251135987

# Row 393
$ws.Cells.Item(393, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(393, 2).Value = 60
$ws.Cells.Item(393, 3).Value = 2900
$ws.Cells.Item(393, 4).Value = 0
$ws.Cells.Item(393, 5).Value = 0.1
$ws.Cells.Item(393, 6).Value = 50.08872597113812
$ws.Cells.Item(393, 7).Value = 1809.243007695047
$ws.Cells.Item(393, 8).Value = 46.05060087590582
$ws.Cells.Item(393, 9).Value = 6.953267749580917
$ws.Cells.Item(393, 10).Value = 10.99139284481321

# Row 394
$ws.Cells.Item(394, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(394, 2).Value = 420
$ws.Cells.Item(394, 3).Value = 2900
$ws.Cells.Item(394, 4).Value = 0
$ws.Cells.Item(394, 5).Value = 0.1
$ws.Cells.Item(394, 6).Value = 48.1874030872552
$ws.Cells.Item(394, 7).Value = 1067.974215221202
$ws.Cells.Item(394, 8).Value = 31.50180904133441
$ws.Cells.Item(394, 9).Value = 8.854590633463836
$ws.Cells.Item(394, 10).Value = 25.54018467938462

# Row 395
$ws.Cells.Item(395, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(395, 2).Value = 3
$ws.Cells.Item(395, 3).Value = 2900
$ws.Cells.Item(395, 4).Value = 0
$ws.Cells.Item(395, 5).Value = 0.1
$ws.Cells.Item(395, 6).Value = 53.87595673098347
$ws.Cells.Item(395, 7).Value = 10385.36486245873
$ws.Cells.Item(395, 8).Value = 52.71698023374482
$ws.Cells.Item(395, 9).Value = 3.166036989735566
$ws.Cells.Item(395, 10).Value = 4.325013486974214

# Row 396
$ws.Cells.Item(396, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(396, 2).Value = 1
$ws.Cells.Item(396, 3).Value = 2900
$ws.Cells.Item(396, 4).Value = 0
$ws.Cells.Item(396, 5).Value = 0.1
$ws.Cells.Item(396, 6).Value = 54.58431692148883
$ws.Cells.Item(396, 7).Value = 21741.68594864994
$ws.Cells.Item(396, 8).Value = 53.77554729363442
$ws.Cells.Item(396, 9).Value = 2.457676799230207
$ws.Cells.Item(396, 10).Value = 3.266446427084617

# Row 397
$ws.Cells.Item(397, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(397, 2).Value = 8
$ws.Cells.Item(397, 3).Value = 2900
$ws.Cells.Item(397, 4).Value = 0
$ws.Cells.Item(397, 5).Value = 0.1
$ws.Cells.Item(397, 6).Value = 53.65701634015197
$ws.Cells.Item(397, 7).Value = 8747.874059747004
$ws.Cells.Item(397, 8).Value = 51.05371690347397
$ws.Cells.Item(397, 9).Value = 3.384977380567065
$ws.Cells.Item(397, 10).Value = 5.988276817245065

# Row 398
$ws.Cells.Item(398, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(398, 2).Value = 6
$ws.Cells.Item(398, 3).Value = 2900
$ws.Cells.Item(398, 4).Value = 0
$ws.Cells.Item(398, 5).Value = 0.1
$ws.Cells.Item(398, 6).Value = 53.71567707197253
$ws.Cells.Item(398, 7).Value = 9009.80646665667
$ws.Cells.Item(398, 8).Value = 51.70474070568839
$ws.Cells.Item(398, 9).Value = 3.326316648746506
$ws.Cells.Item(398, 10).Value = 5.337253015030647

# Row 399
$ws.Cells.Item(399, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(399, 2).Value = 610
$ws.Cells.Item(399, 3).Value = 2900
$ws.Cells.Item(399, 4).Value = 0
$ws.Cells.Item(399, 5).Value = 0.1
$ws.Cells.Item(399, 6).Value = 46.86861350038294
$ws.Cells.Item(399, 7).Value = 1002.245676223203
$ws.Cells.Item(399, 8).Value = 24.12624764429263
$ws.Cells.Item(399, 9).Value = 10.17338022033609
$ws.Cells.Item(399, 10).Value = 32.9157460764264

# Row 400
$ws.Cells.Item(400, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(400, 2).Value = 10000000000000000
$ws.Cells.Item(400, 3).Value = 2900
$ws.Cells.Item(400, 4).Value = 0
$ws.Cells.Item(400, 5).Value = 0.1
$ws.Cells.Item(400, 6).Value = 29.326573568
$ws.Cells.Item(400, 7).Value = 889.6259531912921
$ws.Cells.Item(400, 8).Value = -330932225215385.4
$ws.Cells.Item(400, 9).Value = 27.71542015271903
$ws.Cells.Item(400, 10).Value = 330932225215442.4

# Row 401
$ws.Cells.Item(401, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(401, 2).Value = 630
$ws.Cells.Item(401, 3).Value = 2900
$ws.Cells.Item(401, 4).Value = 0
$ws.Cells.Item(401, 5).Value = 0.1
$ws.Cells.Item(401, 6).Value = 46.70547397112593
$ws.Cells.Item(401, 7).Value = 995.1507456279003
$ws.Cells.Item(401, 8).Value = 23.38372923321111
$ws.Cells.Item(401, 9).Value = 10.3365197495931
$ws.Cells.Item(401, 10).Value = 33.65826448750792

# Row 402
$ws.Cells.Item(402, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(402, 2).Value = 620
$ws.Cells.Item(402, 3).Value = 2900
$ws.Cells.Item(402, 4).Value = 0
$ws.Cells.Item(402, 5).Value = 0.1
$ws.Cells.Item(402, 6).Value = 46.78588689719382
$ws.Cells.Item(402, 7).Value = 998.6204177211238
$ws.Cells.Item(402, 8).Value = 23.75430615319304
$ws.Cells.Item(402, 9).Value = 10.25610682352522
$ws.Cells.Item(402, 10).Value = 33.28768756752599

# Row 403
$ws.Cells.Item(403, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(403, 2).Value = 390
$ws.Cells.Item(403, 3).Value = 2900
$ws.Cells.Item(403, 4).Value = 0
$ws.Cells.Item(403, 5).Value = 0.1
$ws.Cells.Item(403, 6).Value = 48.18740308725521
$ws.Cells.Item(403, 7).Value = 1067.974215221202
$ws.Cells.Item(403, 8).Value = 32.69363718747162
$ws.Cells.Item(403, 9).Value = 8.854590633463822
$ws.Cells.Item(403, 10).Value = 24.34835653324741

# Row 404
$ws.Cells.Item(404, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(404, 2).Value = 4
$ws.Cells.Item(404, 3).Value = 2900
$ws.Cells.Item(404, 4).Value = 0
$ws.Cells.Item(404, 5).Value = 0.1
$ws.Cells.Item(404, 6).Value = 53.80672461550006
$ws.Cells.Item(404, 7).Value = 9519.301348137373
$ws.Cells.Item(404, 8).Value = 52.39028950978646
$ws.Cells.Item(404, 9).Value = 3.235269105218975
$ws.Cells.Item(404, 10).Value = 4.651704210932571

# Row 405
$ws.Cells.Item(405, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(405, 2).Value = 410
$ws.Cells.Item(405, 3).Value = 2900
$ws.Cells.Item(405, 4).Value = 0
$ws.Cells.Item(405, 5).Value = 0.1
$ws.Cells.Item(405, 6).Value = 48.1874030872552
$ws.Cells.Item(405, 7).Value = 1067.974215221202
$ws.Cells.Item(405, 8).Value = 31.89908509004681
$ws.Cells.Item(405, 9).Value = 8.854590633463836
$ws.Cells.Item(405, 10).Value = 25.14290863067222

# Row 406
$ws.Cells.Item(406, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(406, 2).Value = 7
$ws.Cells.Item(406, 3).Value = 2900
$ws.Cells.Item(406, 4).Value = 0
$ws.Cells.Item(406, 5).Value = 0.1
$ws.Cells.Item(406, 6).Value = 53.49831041749932
$ws.Cells.Item(406, 7).Value = 8717.176637523147
$ws.Cells.Item(406, 8).Value = 51.22841681071766
$ws.Cells.Item(406, 9).Value = 3.543683303219716
$ws.Cells.Item(406, 10).Value = 5.813576910001373

# Row 407
$ws.Cells.Item(407, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(407, 2).Value = 145
$ws.Cells.Item(407, 3).Value = 2900
$ws.Cells.Item(407, 4).Value = 0
$ws.Cells.Item(407, 5).Value = 0.1
$ws.Cells.Item(407, 6).Value = 48.18740308725521
$ws.Cells.Item(407, 7).Value = 1067.974215221202
$ws.Cells.Item(407, 8).Value = 42.42690038092541
$ws.Cells.Item(407, 9).Value = 8.854590633463822
$ws.Cells.Item(407, 10).Value = 14.61509333979362

# Row 408
$ws.Cells.Item(408, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(408, 2).Value = 740
$ws.Cells.Item(408, 3).Value = 2900
$ws.Cells.Item(408, 4).Value = 0
$ws.Cells.Item(408, 5).Value = 0.1
$ws.Cells.Item(408, 6).Value = 45.88276327422266
$ws.Cells.Item(408, 7).Value = 962.6327851215189
$ws.Cells.Item(408, 8).Value = 19.38409894321065
$ws.Cells.Item(408, 9).Value = 11.15923044649637
$ws.Cells.Item(408, 10).Value = 37.65789477750839

# Row 409
$ws.Cells.Item(409, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(409, 2).Value = 370
$ws.Cells.Item(409, 3).Value = 2900
$ws.Cells.Item(409, 4).Value = 0
$ws.Cells.Item(409, 5).Value = 0.1
$ws.Cells.Item(409, 6).Value = 48.18740308725521
$ws.Cells.Item(409, 7).Value = 1067.974215221202
$ws.Cells.Item(409, 8).Value = 33.48818928489641
$ws.Cells.Item(409, 9).Value = 8.854590633463822
$ws.Cells.Item(409, 10).Value = 23.55380443582263

# Row 410
$ws.Cells.Item(410, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(410, 2).Value = 600
$ws.Cells.Item(410, 3).Value = 2900
$ws.Cells.Item(410, 4).Value = 0
$ws.Cells.Item(410, 5).Value = 0.1
$ws.Cells.Item(410, 6).Value = 46.95379500295929
$ws.Cells.Item(410, 7).Value = 1006.038007987869
$ws.Cells.Item(410, 8).Value = 24.4996121787907
$ws.Cells.Item(410, 9).Value = 10.08819871775975
$ws.Cells.Item(410, 10).Value = 32.54238154192834

# Row 411
$ws.Cells.Item(411, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(411, 2).Value = 680
$ws.Cells.Item(411, 3).Value = 2900
$ws.Cells.Item(411, 4).Value = 0
$ws.Cells.Item(411, 5).Value = 0.1
$ws.Cells.Item(411, 6).Value = 46.35242482763016
$ws.Cells.Item(411, 7).Value = 980.5402281018844
$ws.Cells.Item(411, 8).Value = 21.54932598189934
$ws.Cells.Item(411, 9).Value = 10.68956889308888
$ws.Cells.Item(411, 10).Value = 35.49266773881969

# Row 412
$ws.Cells.Item(412, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(412, 2).Value = 430
$ws.Cells.Item(412, 3).Value = 2900
$ws.Cells.Item(412, 4).Value = 0
$ws.Cells.Item(412, 5).Value = 0.1
$ws.Cells.Item(412, 6).Value = 48.18740308725521
$ws.Cells.Item(412, 7).Value = 1067.974215221202
$ws.Cells.Item(412, 8).Value = 31.10453299262203
$ws.Cells.Item(412, 9).Value = 8.854590633463822
$ws.Cells.Item(412, 10).Value = 25.937460728097

# Row 413
$ws.Cells.Item(413, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(413, 2).Value = 730
$ws.Cells.Item(413, 3).Value = 2900
$ws.Cells.Item(413, 4).Value = 0
$ws.Cells.Item(413, 5).Value = 0.1
$ws.Cells.Item(413, 6).Value = 45.95906874431417
$ws.Cells.Item(413, 7).Value = 965.4264344845503
$ws.Cells.Item(413, 8).Value = 19.74263210790998
$ws.Cells.Item(413, 9).Value = 11.08292497640486
$ws.Cells.Item(413, 10).Value = 37.29936161280905

# Row 414
$ws.Cells.Item(414, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(414, 2).Value = 165
$ws.Cells.Item(414, 3).Value = 2900
$ws.Cells.Item(414, 4).Value = 0
$ws.Cells.Item(414, 5).Value = 0.1
$ws.Cells.Item(414, 6).Value = 48.18740308725522
$ws.Cells.Item(414, 7).Value = 1067.974215221202
$ws.Cells.Item(414, 8).Value = 41.63234828350063
$ws.Cells.Item(414, 9).Value = 8.854590633463815
$ws.Cells.Item(414, 10).Value = 15.4096454372184

# Row 415
$ws.Cells.Item(415, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(415, 2).Value = 400
$ws.Cells.Item(415, 3).Value = 2900
$ws.Cells.Item(415, 4).Value = 0
$ws.Cells.Item(415, 5).Value = 0.1
$ws.Cells.Item(415, 6).Value = 48.18740308725521
$ws.Cells.Item(415, 7).Value = 1067.974215221202
$ws.Cells.Item(415, 8).Value = 32.29636113875923
$ws.Cells.Item(415, 9).Value = 8.854590633463822
$ws.Cells.Item(415, 10).Value = 24.7456325819598

# Row 416
$ws.Cells.Item(416, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(416, 2).Value = 640
$ws.Cells.Item(416, 3).Value = 2900
$ws.Cells.Item(416, 4).Value = 0
$ws.Cells.Item(416, 5).Value = 0.1
$ws.Cells.Item(416, 6).Value = 46.62724313095844
$ws.Cells.Item(416, 7).Value = 991.8261512019781
$ws.Cells.Item(416, 8).Value = 23.01446184881642
$ws.Cells.Item(416, 9).Value = 10.41475058976059
$ws.Cells.Item(416, 10).Value = 34.02753187190261

# Row 417
$ws.Cells.Item(417, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(417, 2).Value = 75
$ws.Cells.Item(417, 3).Value = 2900
$ws.Cells.Item(417, 4).Value = 0
$ws.Cells.Item(417, 5).Value = 0.1
$ws.Cells.Item(417, 6).Value = 48.1874030872552
$ws.Cells.Item(417, 7).Value = 1067.974215221202
$ws.Cells.Item(417, 8).Value = 45.2078327219122
$ws.Cells.Item(417, 9).Value = 8.854590633463836
$ws.Cells.Item(417, 10).Value = 11.83416099880684

# Row 418
$ws.Cells.Item(418, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(418, 2).Value = 55
$ws.Cells.Item(418, 3).Value = 2900
$ws.Cells.Item(418, 4).Value = 0
$ws.Cells.Item(418, 5).Value = 0.1
$ws.Cells.Item(418, 6).Value = 50.59192765071704
$ws.Cells.Item(418, 7).Value = 2044.961440343802
$ws.Cells.Item(418, 8).Value = 46.40804564246646
$ws.Cells.Item(418, 9).Value = 6.450066070001995
$ws.Cells.Item(418, 10).Value = 10.63394807825257

# Row 419
$ws.Cells.Item(419, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(419, 2).Value = 125
$ws.Cells.Item(419, 3).Value = 2900
$ws.Cells.Item(419, 4).Value = 0
$ws.Cells.Item(419, 5).Value = 0.1
$ws.Cells.Item(419, 6).Value = 48.18740308725521
$ws.Cells.Item(419, 7).Value = 1067.974215221202
$ws.Cells.Item(419, 8).Value = 43.22145247835022
$ws.Cells.Item(419, 9).Value = 8.854590633463822
$ws.Cells.Item(419, 10).Value = 13.82054124236881

# Row 420
$ws.Cells.Item(420, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(420, 2).Value = 105
$ws.Cells.Item(420, 3).Value = 2900
$ws.Cells.Item(420, 4).Value = 0
$ws.Cells.Item(420, 5).Value = 0.1
$ws.Cells.Item(420, 6).Value = 48.1874030872552
$ws.Cells.Item(420, 7).Value = 1067.974215221202
$ws.Cells.Item(420, 8).Value = 44.016004575775
$ws.Cells.Item(420, 9).Value = 8.854590633463836
$ws.Cells.Item(420, 10).Value = 13.02598914494403

# Row 421
$ws.Cells.Item(421, 1).Value = "linearization_heuristic_optBouncing=False_initial_uhat=time_gradient_targetGroups=False_targetAct=False_targetTests=True"
$ws.Cells.Item(421, 2).Value = 9
$ws.Cells.Item(421, 3).Value = 2900
$ws.Cells.Item(421, 4).Value = 0
$ws.Cells.Item(421, 5).Value = 0.1
$ws.Cells.Item(421, 6).Value = 53.61810498285629
$ws.Cells.Item(421, 7).Value = 8596.69973313073
$ws.Cells.Item(421, 8).Value = 50.74000496139278
$ws.Cells.Item(421, 9).Value = 3.423888737862747
$ws.Cells.Item(421, 10).Value = 6.301988759326257
